$wb = $excel.ActiveWorkbook
$example = $wb.Worksheets.Item("Example")
$count = $wb.Worksheets.Item("Count")

# --- Example sheet: update Sprint # 2 row (row 5) ---
$example.Range("D5").Value = "A, B, C, E, F, H, I, J, L, M, N"
$example.Range("E5").Value = "Configuring mongodb/cloud server, deploying app to cloud platform, user roles and auth, training users on API, installed additional security protocols"
$example.Range("G5").Value = "github.com/mblaul/skypi/api,`r`nGoogle Team Drive,`r`nPostman Documentation`r`n"
$example.Range("H5").Value = "/core/api/*,`r`n/core/api/models/User.js,`r`n/core/api/models/Weather.js,`r`n/core/api/controllers/user.js,`r`n/core/api/controllers/weather.js"

# --- Count sheet: update sprint tallies ---
$count.Range("A5").Value = 1

$count.Range("A6").Value = 2
$count.Range("C6").Value = 1
$count.Range("D6").Value = 1
$count.Range("E6").Value = 3
$count.Range("G6").Value = 1
$count.Range("H6").Value = 2
$count.Range("J6").Value = 1
$count.Range("K6").Value = 6
$count.Range("L6").Value = 2
$count.Range("O6").Value = 1
$count.Range("P6").Value = 1

$count.Range("A7").Value = 3
$count.Range("A8").Value = 4
$count.Range("A9").Value = 5
$count.Range("A10").Value = 6
$count.Range("A11").Value = 7
$count.Range("A12").Value = 8
$count.Range("A13").Value = 9

# --- Switch active tab from Count back to Example, update selections ---
$count.Range("D9").Select()
$example.Activate()
$example.Range("I5").Select()
